$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Debug" config row (row 2): start_x (Start Time), Acceleration Time,
# plateau (Plateau Time), Deceleration Time, High Level Frequency and end_x (End Time)
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 40
$ws.Range("J2").Value = 3

# Update view/selection state to match the saved workbook
$ws.Range("J3").Select()
